$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '29.023.28'
Set-TextValue $ws 'E2' '  -0.69%  '
Set-TextValue $ws 'D3' '1.817.00'
Set-TextValue $ws 'E3' '  -0.73%  '
Set-TextValue $ws 'E4' '  +0.33%  '
Set-TextValue $ws 'D5' '232.83'
Set-TextValue $ws 'E5' '  -1.86%  '
Set-TextValue $ws 'D6' '0.5902'
Set-TextValue $ws 'E6' '  -3.15%  '
Set-TextValue $ws 'D7' '1.004'
Set-TextValue $ws 'E7' '  +0.25%  '
Set-TextValue $ws 'D8' '0.2745'
Set-TextValue $ws 'E8' '  -2.58%  '
Set-TextValue $ws 'D9' '0.06761'
Set-TextValue $ws 'E9' '  -4.66%  '
Set-TextValue $ws 'D10' '22.92'
Set-TextValue $ws 'E10' '  -3.96%  '
Set-TextValue $ws 'D11' '0.07491'
Set-TextValue $ws 'E11' '  -2.03%  '
Set-TextValue $ws 'D12' '1.816.41'
Set-TextValue $ws 'E12' '  -1.13%  '
Set-TextValue $ws 'D13' '4.669'
Set-TextValue $ws 'E13' '  -2.86%  '
Set-TextValue $ws 'D14' '0.6216'
Set-TextValue $ws 'E14' '  -1.92%  '
Set-TextValue $ws 'D15' '0.000009386'
Set-TextValue $ws 'E15' '  -6.16%  '
Set-TextValue $ws 'D16' '74.41'
Set-TextValue $ws 'E16' '  -6.50%  '
Set-TextValue $ws 'D17' '28.740.78'
Set-TextValue $ws 'E17' '  -1.62%  '
Set-TextValue $ws 'D18' '5.421'
Set-TextValue $ws 'E18' '  -9.03%  '
Set-TextValue $ws 'D19' '1.004'
Set-TextValue $ws 'E19' '  +0.25%  '
Set-TextValue $ws 'D20' '207.37'
Set-TextValue $ws 'E20' '  -9.44%  '
Set-TextValue $ws 'D21' '11.35'
Set-TextValue $ws 'E21' '  -3.89%  '
Set-TextValue $ws 'D22' '6.755'
Set-TextValue $ws 'E22' '  -3.98%  '
Set-TextValue $ws 'E23' '  +0.42%  '
Set-TextValue $ws 'D24' '155.05'
Set-TextValue $ws 'E24' '  -0.23%  '
Set-TextValue $ws 'D25' '0.1267'
Set-TextValue $ws 'E25' '  -2.31%  '
Set-TextValue $ws 'D26' '7.761'
Set-TextValue $ws 'E26' '  -4.20%  '
Set-TextValue $ws 'D27' '16.24'
Set-TextValue $ws 'E27' '  -2.85%  '
Set-TextValue $ws 'D28' '0.06438'
Set-TextValue $ws 'E28' '  -5.38%  '
Set-TextValue $ws 'D29' '1.400'
Set-TextValue $ws 'E29' '  -5.30%  '
Set-TextValue $ws 'D30' '1.430'
Set-TextValue $ws 'E30' '  -1.95%  '
Set-TextValue $ws 'D31' '3.708'
Set-TextValue $ws 'E31' '  -3.15%  '
Set-TextValue $ws 'D32' '3.668'
Set-TextValue $ws 'E32' '  -4.19%  '
Set-TextValue $ws 'D33' '1.676'
Set-TextValue $ws 'E33' '  -2.63%  '
Set-TextValue $ws 'D34' '1.046'
Set-TextValue $ws 'E34' '  -7.22%  '
Set-TextValue $ws 'D36' '0.6280'
Set-TextValue $ws 'E36' '  -4.21%  '
Set-TextValue $ws 'D37' '2.740'
Set-TextValue $ws 'E37' '  -0.91%  '
Set-TextValue $ws 'D38' '6.413'
Set-TextValue $ws 'E38' '  -2.70%  '
Set-TextValue $ws 'D39' '0.01696'
Set-TextValue $ws 'E39' '  -4.02%  '
Set-TextValue $ws 'D40' '1.127.20'
Set-TextValue $ws 'E40' '  -8.61%  '
Set-TextValue $ws 'D41' '0.8678'
Set-TextValue $ws 'E41' '  -6.05%  '
Set-TextValue $ws 'D42' '1.004'
Set-TextValue $ws 'E42' '  +0.24%  '
Set-TextValue $ws 'D43' '1.972.29'
Set-TextValue $ws 'E43' '  -0.53%  '
Set-TextValue $ws 'D44' '99.98'
Set-TextValue $ws 'E44' '  -0.83%  '
Set-TextValue $ws 'D45' '60.05'
Set-TextValue $ws 'E45' '  -5.38%  '
Set-TextValue $ws 'E46' '  -3.92%  '
Set-TextValue $ws 'D47' '1.567'
Set-TextValue $ws 'E47' '  -3.77%  '
Set-TextValue $ws 'D48' '0.05465'
Set-TextValue $ws 'E48' '  -1.72%  '
Set-TextValue $ws 'D49' '0.4509'
Set-TextValue $ws 'E49' '  -1.23%  '
Set-TextValue $ws 'E50' '  -3.43%  '
Set-TextValue $ws 'D51' '1.002'
Set-TextValue $ws 'E51' '  +0.32%  '
